# TD-4469 Limits details query by learner centre ID and updates related tests
# Rename the "Usage Statistics" sheet to "Usage summary" to match the
# updated workbook.xml <sheet> entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Usage summary"
